# Update ICM to v4, add '96 consists
# - Re-sort the "96-97" sheet by column A (Materieel) instead of column B (Treinserie)
# - Make "96-97" the active sheet/tab with the selection on A15
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("96-97")
$ws.Activate()

$rng = $ws.Range("A1:C77")
$rng.Sort($ws.Range("A1"), 1, $null, $null, 1, $null, 1, 1)

$ws.Range("A15").Select()
